# Re-upload of the full program: rows got reshuffled in two sheets of the
# workbook (CLASIFICACION_ABC+D_TIERRA_ARIDOS_P2_2025.xlsx) during the
# re-export. This script rewrites the affected rows so the data matches the
# new row order, cell by cell.

$wb = $excel.ActiveWorkbook

function Set-RowData($sheet, $rowNum, $data) {
    foreach ($col in $data.Keys) {
        $addr = $col + $rowNum
        $sheet.Range($addr).Value = $data[$col]
    }
}

# ---------------------------------------------------------------------
# Sheet "CATEGORIA C - BAJO IMPACTO" (3rd sheet): rows 46 and 47 swap
# their data.
# ---------------------------------------------------------------------
$wsC = $wb.Worksheets.Item(3)

$rowC46 = @{
    A = "3203050026"
    B = "SACO MARMOLINA AMARILLO 20KG"
    C = "9I12"
    D = "UNICO"
    E = "32"
    F = "MANTENIMIENTO"
    G = 90
    H = 1
    I = 5.99
    J = 2.48
    K = 0
    L = 0
    M = 0.5
    N = 1.5
    O = -1
    P = 44
    Q = 0
    R = 0
    S = 0
    T = "Bajo"
    U = "AUMENTAR STOCK: Producto de alto interés. Incrementar compras 30% próxima temporada. Stock actual: -1 unidades. Stock objetivo: 1 unidades. Alta rotación confirmada."
    V = "Sin stock"
    W = "25"
}

$rowC47 = @{
    A = "3101010008"
    B = "SUSTRATO UNIVERSAL ECOLOGICO SIN TURBA"
    C = "10L"
    D = "UNICO"
    E = "31"
    F = "TIERRAS"
    G = 90
    H = 2
    I = 6.58
    J = 2.96
    K = 0
    L = 0
    M = 1
    N = 2.9
    O = -2
    P = 55
    Q = 0
    R = 0
    S = 0
    T = "Bajo"
    U = "AUMENTAR STOCK: Producto de alto interés. Incrementar compras 30% próxima temporada. Stock actual: -2 unidades. Stock objetivo: 1 unidades. Alta rotación confirmada."
    V = "Sin stock"
    W = "25"
}

Set-RowData $wsC 46 $rowC46
Set-RowData $wsC 47 $rowC47

# ---------------------------------------------------------------------
# Sheet "CATEGORIA D - SIN VENTAS" (4th sheet): rows 4-10 are reshuffled
# into a new order.
# ---------------------------------------------------------------------
$wsD = $wb.Worksheets.Item(4)

$rowD4 = @{
    A = "3102110006"
    B = "BIG BAG TIERRA ENRIQUECIDA 500L (NO VENTA)"
    C = ""
    D = ""
    E = "31"
    F = "TIERRAS"
    G = 90
    H = 0
    I = 0
    J = 0
    K = 0
    L = 7
    M = 0
    N = 0
    O = 7
    P = 92
    Q = 12
    R = 13.33
    S = 0
    T = "Crítico"
    U = "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
    V = "Compra 19/05/2025"
    W = "14"
}

$rowD5 = @{
    A = "3203050005"
    B = "BIG BAG GRAVA VOLCANICA MARRON 500L (NO VENTA)"
    C = "5I10"
    D = "UNICO"
    E = "32"
    F = "MANTENIMIENTO"
    G = 90
    H = 0
    I = 0
    J = 0
    K = 0
    L = 3
    M = 0
    N = 0
    O = 3
    P = 92
    Q = 12
    R = 13.33
    S = 0
    T = "Crítico"
    U = "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
    V = "Compra 19/05/2025"
    W = "14"
}

$rowD6 = @{
    A = "3201020001"
    B = "BIG BAG MANTILLO 500L (NO VENTA)"
    C = ""
    D = ""
    E = "32"
    F = "MANTENIMIENTO"
    G = 90
    H = 0
    I = 0
    J = 0
    K = 0
    L = 7
    M = 0
    N = 0
    O = 7
    P = 92
    Q = 12
    R = 13.33
    S = 0
    T = "Crítico"
    U = "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
    V = "Compra 19/05/2025"
    W = "14"
}

$rowD7 = @{
    A = "3101010010"
    B = "COMPO BIO SUBSTRATO AQUA DEPOT"
    C = ""
    D = ""
    E = "31"
    F = "TIERRAS"
    G = 90
    H = 0
    I = 0
    J = 0
    K = 0
    L = 41
    M = 0
    N = 0
    O = 41
    P = 92
    Q = 92
    R = 102.22
    S = 20
    T = "Crítico"
    U = "LIQUIDACIÓN URGENTE: Aplicar descuento 20% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 244.24€. Prioridad máxima."
    V = "Stock inicial"
    W = "14"
}

$rowD8 = @{
    A = "3102110006"
    B = "BIG BAG TIERRA ENRIQUECIDA 650L (NO VENTA)"
    C = ""
    D = ""
    E = "31"
    F = "TIERRAS"
    G = 90
    H = 0
    I = 0
    J = 0
    K = 0
    L = 10
    M = 0
    N = 0
    O = 10
    P = 92
    Q = 92
    R = 102.22
    S = 20
    T = "Crítico"
    U = "LIQUIDACIÓN URGENTE: Aplicar descuento 20% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 144.76€. Prioridad máxima."
    V = "Stock inicial"
    W = "14"
}

$rowD9 = @{
    A = "3102110001"
    B = "BIG BAG JABRE"
    C = "1M3"
    D = "UNICO"
    E = "31"
    F = "TIERRAS"
    G = 90
    H = 0
    I = 0
    J = 0
    K = 0
    L = 14
    M = 0
    N = 0
    O = 14
    P = 92
    Q = 92
    R = 102.22
    S = 20
    T = "Crítico"
    U = "LIQUIDACIÓN URGENTE: Aplicar descuento 20% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 458.74€. Prioridad máxima."
    V = "Stock inicial"
    W = "14"
}

$rowD10 = @{
    A = "3202030010"
    B = "BIG BAG RECEBO CESPED 650L (NO VENTA)"
    C = ""
    D = ""
    E = "32"
    F = "MANTENIMIENTO"
    G = 90
    H = 0
    I = 0
    J = 0
    K = 0
    L = 1
    M = 0
    N = 0
    O = 1
    P = 92
    Q = 92
    R = 102.22
    S = 20
    T = "Crítico"
    U = "LIQUIDACIÓN URGENTE: Aplicar descuento 20% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 16.0€. Prioridad máxima."
    V = "Stock inicial"
    W = "14"
}

Set-RowData $wsD 4 $rowD4
Set-RowData $wsD 5 $rowD5
Set-RowData $wsD 6 $rowD6
Set-RowData $wsD 7 $rowD7
Set-RowData $wsD 8 $rowD8
Set-RowData $wsD 9 $rowD9
Set-RowData $wsD 10 $rowD10
